# Update "want to go" counts (column F) on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$sheetExhibition = $wb.Worksheets.Item("展览")
$sheetAllTypes   = $wb.Worksheets.Item("全部类型")

# Row => new F value, for sheet "展览"
$exhibitionUpdates = @{
    2  = 92
    3  = 4043
    4  = 2366
    5  = 471
    10 = 119
    11 = 79
    12 = 133
    13 = 1510
    15 = 2866
}

foreach ($row in $exhibitionUpdates.Keys) {
    $sheetExhibition.Range("F$row").Value = $exhibitionUpdates[$row]
}

# Row => new F value, for sheet "全部类型"
$allTypesUpdates = @{
    2  = 92
    3  = 4043
    4  = 2366
    5  = 471
    11 = 119
    12 = 79
    13 = 133
    16 = 1510
    18 = 2866
}

foreach ($row in $allTypesUpdates.Keys) {
    $sheetAllTypes.Range("F$row").Value = $allTypesUpdates[$row]
}
